# Weekly update: insert a new price record as row 24 ("Asterix" / "1a (cosecha
# lavada)" quote for a later date), pushing the existing rows 24-57 down to
# rows 25-58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 24; Excel shifts rows
# 24..57 down to 25..58 and carries the row's number formatting along
# (column D keeps its date style).
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new data point.
$ws.Cells.Item(24, 1).Value  = 1
$ws.Cells.Item(24, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(24, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(24, 4).Value  = 44571
$ws.Cells.Item(24, 5).Value  = 15
$ws.Cells.Item(24, 6).Value  = 100114001
$ws.Cells.Item(24, 7).Value  = "Papa"
$ws.Cells.Item(24, 8).Value  = "Asterix"
$ws.Cells.Item(24, 9).Value  = "1a (cosecha lavada)"
$ws.Cells.Item(24, 10).Value = 1000
$ws.Cells.Item(24, 11).Value = 14000
$ws.Cells.Item(24, 12).Value = 15000
$ws.Cells.Item(24, 13).Value = 14500
$ws.Cells.Item(24, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(24, 16).Value = 580
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
